# Applies the "Terminologia_glossario" ID fix: decrements several
# ID_db_Id_statico_entry (column C) values by 1 for rows 35-65, and
# ID_db_Commento_entry (column B) values by 1 for rows 59-65.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (ID_db_Id_statico_entry) updates for rows 35-65
$ws.Range("C35").Value = 14000039
$ws.Range("C36").Value = 14000040
$ws.Range("C37").Value = 14000043
$ws.Range("C38").Value = 14000044
$ws.Range("C39").Value = 14000045
$ws.Range("C40").Value = 14000046
$ws.Range("C41").Value = 14000051
$ws.Range("C42").Value = 14000054
$ws.Range("C43").Value = 14000055
$ws.Range("C44").Value = 14000056
$ws.Range("C45").Value = 14000057
$ws.Range("C46").Value = 14000058
$ws.Range("C47").Value = 14000059
$ws.Range("C48").Value = 14000060
$ws.Range("C49").Value = 14000061
$ws.Range("C50").Value = 14000062
$ws.Range("C51").Value = 14000063
$ws.Range("C52").Value = 14000064
$ws.Range("C53").Value = 14000065
$ws.Range("C54").Value = 14000068
$ws.Range("C55").Value = 14000071
$ws.Range("C56").Value = 14000073
$ws.Range("C57").Value = 14000076
$ws.Range("C58").Value = 14000074
$ws.Range("C59").Value = 14000041
$ws.Range("C60").Value = 14000047
$ws.Range("C61").Value = 14000049
$ws.Range("C62").Value = 14000052
$ws.Range("C63").Value = 14000066
$ws.Range("C64").Value = 14000069
$ws.Range("C65").Value = 14000075

# Column B (ID_db_Commento_entry) updates for rows 59-65
$ws.Range("B59").Value = 12000041
$ws.Range("B60").Value = 12000047
$ws.Range("B61").Value = 12000049
$ws.Range("B62").Value = 12000052
$ws.Range("B63").Value = 12000066
$ws.Range("B64").Value = 12000069
$ws.Range("B65").Value = 12000075
